# codeforIATI SectorGroup.xlsx codelist update:
# The columns E (group-name), F (category-name) and G (group-code) are
# rotated on every row (including the header) so that:
#   new E = old G   (group-code moves into the E column)
#   new F = old E   (group-name moves into the F column)
#   new G = old F   (category-name moves into the G column)
#
# Using Range.Copy (instead of plain Value assignment) preserves each
# cell's original text/number type so numeric-looking strings like "110"
# stay stored as text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = $ws.UsedRange.Row
$lastRow = $firstRow + $ws.UsedRange.Rows.Count - 1

# Scratch cell well outside the used range (A1:G235) used as temporary
# holding space while rotating E/F/G for a row.
$scratch = $ws.Cells.Item(1, 26)

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $gCell.Copy($scratch)   # scratch = old G
    $fCell.Copy($gCell)     # new G = old F
    $eCell.Copy($fCell)     # new F = old E
    $scratch.Copy($eCell)   # new E = old G
}

$scratch.ClearContents()
